$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object 'object[,]' 50,4
$values[0,0] = 0.0891108363866806
$values[0,1] = 0.9776841402053833
$values[0,2] = 0.2068810611963272
$values[0,3] = 0.9546296000480652
$values[1,0] = 0.02537666819989681
$values[1,1] = 0.9900674819946289
$values[1,2] = 0.1479210704565048
$values[1,3] = 0.949579119682312
$values[2,0] = 0.01584483310580254
$values[2,1] = 0.9947757720947266
$values[2,2] = 0.09521707147359848
$values[2,3] = 0.9622054100036621
$values[3,0] = 0.01257229223847389
$values[3,1] = 0.9955711960792542
$values[3,2] = 0.354936808347702
$values[3,3] = 0.9484006762504578
$values[4,0] = 0.01130461040884256
$values[4,1] = 0.995549738407135
$values[4,2] = 0.1660708338022232
$values[4,3] = 0.9659090638160706
$values[5,0] = 0.01128542423248291
$values[5,1] = 0.9952917695045471
$values[5,2] = 0.5267589688301086
$values[5,3] = 0.9485689997673035
$values[6,0] = 0.009384628385305405
$values[6,1] = 0.9959797263145447
$values[6,2] = 0.4426897168159485
$values[6,3] = 0.9493265748023987
$values[7,0] = 0.009427289478480816
$values[7,1] = 0.9953347444534302
$values[7,2] = 0.4434752464294434
$values[7,3] = 0.9488215446472168
$values[8,0] = 0.008229215629398823
$values[8,1] = 0.9954422116279602
$values[8,2] = 0.3345008790493011
$values[8,3] = 0.9666666388511658
$values[9,0] = 0.007301429286599159
$values[9,1] = 0.9957432150840759
$values[9,2] = 0.3620936274528503
$values[9,3] = 0.957239031791687
$values[10,0] = 0.007760378532111645
$values[10,1] = 0.9955927133560181
$values[10,2] = 0.1539175063371658
$values[10,3] = 0.9722222089767456
$values[11,0] = 0.007689092773944139
$values[11,1] = 0.995549738407135
$values[11,2] = 0.2370198369026184
$values[11,3] = 0.9748316407203674
$values[12,0] = 0.007159396074712276
$values[12,1] = 0.995614230632782
$values[12,2] = 0.2719425559043884
$values[12,3] = 0.9626262784004211
$values[13,0] = 0.008134615607559681
$values[13,1] = 0.995721697807312
$values[13,2] = 0.246554434299469
$values[13,3] = 0.9588384032249451
$values[14,0] = 0.006642026361078024
$values[14,1] = 0.9959582090377808
$values[14,2] = 0.249793753027916
$values[14,3] = 0.959764301776886
$values[15,0] = 0.006941849365830421
$values[15,1] = 0.9955927133560181
$values[15,2] = 0.3751275241374969
$values[15,3] = 0.9514309763908386
$values[16,0] = 0.007286000531166792
$values[16,1] = 0.9960871934890747
$values[16,2] = 0.1239943578839302
$values[16,3] = 0.9731481671333313
$values[17,0] = 0.00691929180175066
$values[17,1] = 0.9960011839866638
$values[17,2] = 0.1467976719141006
$values[17,3] = 0.9691919088363647
$values[18,0] = 0.007032891735434532
$values[18,1] = 0.9960011839866638
$values[18,2] = 0.1006443798542023
$values[18,3] = 0.9794612526893616
$values[19,0] = 0.007932950742542744
$values[19,1] = 0.9955067038536072
$values[19,2] = 0.6322796940803528
$values[19,3] = 0.9505050778388977
$values[20,0] = 0.006929574068635702
$values[20,1] = 0.995721697807312
$values[20,2] = 0.3654822707176208
$values[20,3] = 0.9677609205245972
$values[21,0] = 0.006169342901557684
$values[21,1] = 0.9959366917610168
$values[21,2] = 0.3533851802349091
$values[21,3] = 0.9682660102844238
$values[22,0] = 0.006651843432337046
$values[22,1] = 0.9960011839866638
$values[22,2] = 0.4587613940238953
$values[22,3] = 0.9558922648429871
$values[23,0] = 0.006962025072425604
$values[23,1] = 0.995485246181488
$values[23,2] = 0.824250340461731
$values[23,3] = 0.949158251285553
$values[24,0] = 0.006739129312336445
$values[24,1] = 0.9962376952171326
$values[24,2] = 0.3952098786830902
$values[24,3] = 0.9506734013557434
$values[25,0] = 0.007398216985166073
$values[25,1] = 0.9952272176742554
$values[25,2] = 0.3174797594547272
$values[25,3] = 0.9651514887809753
$values[26,0] = 0.006100672762840986
$values[26,1] = 0.9959797263145447
$values[26,2] = 0.2634152472019196
$values[26,3] = 0.9724747538566589
$values[27,0] = 0.006198249757289886
$values[27,1] = 0.9965816736221313
$values[27,2] = 0.286062479019165
$values[27,3] = 0.9700336456298828
$values[28,0] = 0.007028468418866396
$values[28,1] = 0.9960011839866638
$values[28,2] = 0.6281360387802124
$values[28,3] = 0.949999988079071
$values[29,0] = 0.00585471885278821
$values[29,1] = 0.9963666796684265
$values[29,2] = 0.6082448959350586
$values[29,3] = 0.9496632814407349
$values[30,0] = 0.007108170073479414
$values[30,1] = 0.995850682258606
$values[30,2] = 0.6204532980918884
$values[30,3] = 0.9520202279090881
$values[31,0] = 0.005757453851401806
$values[31,1] = 0.9963021874427795
$values[31,2] = 0.8222138285636902
$values[31,3] = 0.9494949579238892
$values[32,0] = 0.006630108691751957
$values[32,1] = 0.9958077073097229
$values[32,2] = 0.6453776955604553
$values[32,3] = 0.949158251285553
$values[33,0] = 0.005636299028992653
$values[33,1] = 0.9960656762123108
$values[33,2] = 0.404699832201004
$values[33,3] = 0.9655724167823792
$values[34,0] = 0.006809841375797987
$values[34,1] = 0.995721697807312
$values[34,2] = 0.5373790264129639
$values[34,3] = 0.9499158263206482
$values[35,0] = 0.006198599468916655
$values[35,1] = 0.995850682258606
$values[35,2] = 0.4092240035533905
$values[35,3] = 0.9635521769523621
$values[36,0] = 0.006725901737809181
$values[36,1] = 0.9959797263145447
$values[36,2] = 0.512235701084137
$values[36,3] = 0.9558922648429871
$values[37,0] = 0.006455939263105392
$values[37,1] = 0.9959152340888977
$values[37,2] = 0.6489686965942383
$values[37,3] = 0.9515151381492615
$values[38,0] = 0.006081179715692997
$values[38,1] = 0.995850682258606
$values[38,2] = 0.4167316854000092
$values[38,3] = 0.9673400521278381
$values[39,0] = 0.006706756539642811
$values[39,1] = 0.9961516857147217
$values[39,2] = 0.5394155383110046
$values[39,3] = 0.949158251285553
$values[40,0] = 0.007295573595911264
$values[40,1] = 0.9952057600021362
$values[40,2] = 1.406392216682434
$values[40,3] = 0.9515992999076843
$values[41,0] = 0.006539024412631989
$values[41,1] = 0.995850682258606
$values[41,2] = 1.176237940788269
$values[41,3] = 0.95387202501297
$values[42,0] = 0.005816465243697166
$values[42,1] = 0.9964096546173096
$values[42,2] = 1.369919300079346
$values[42,3] = 0.9505892395973206
$values[43,0] = 0.006660385057330132
$values[43,1] = 0.995657205581665
$values[43,2] = 1.361505031585693
$values[43,3] = 0.9498316645622253
$values[44,0] = 0.007063128054141998
$values[44,1] = 0.995721697807312
$values[44,2] = 1.026727914810181
$values[44,3] = 0.9621211886405945
$values[45,0] = 0.006223268806934357
$values[45,1] = 0.9955711960792542
$values[45,2] = 1.167361378669739
$values[45,3] = 0.9585016965866089
$values[46,0] = 0.006413535214960575
$values[46,1] = 0.9960871934890747
$values[46,2] = 1.173216342926025
$values[46,3] = 0.9577441215515137
$values[47,0] = 0.00701660243794322
$values[47,1] = 0.9955067038536072
$values[47,2] = 1.109143972396851
$values[47,3] = 0.9526935815811157
$values[48,0] = 0.007396447006613016
$values[48,1] = 0.995549738407135
$values[48,2] = 1.119286417961121
$values[48,3] = 0.9515151381492615
$values[49,0] = 0.006853641476482153
$values[49,1] = 0.9960656762123108
$values[49,2] = 0.9749740362167358
$values[49,3] = 0.9508417248725891

$range = $ws.Range("A2:D51")
$range.Value = $values
